# Detection template update: add ClinEpi/ontology parent-term columns, rework
# header/comment rows, and change the sample type example from "stool" to "blood".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 1 (headers)
# ---------------------------------------------------------------------------
# Columns A-H are unchanged. I1 ("website label") is unchanged.
# K1/L1 are overwritten with new header text, J1/M1/N1 are brand-new headers.
$ws.Range("J1").Value = "website parent"
$ws.Range("K1").Value = "website grandparent"
$ws.Range("L1").Value = "<--ClinEpi output | ontology output-->"
$ws.Range("M1").Value = "ontological label"
$ws.Range("N1").Value = "ontological definition"

$ws.Range("J1:N1").WrapText = $true

# ---------------------------------------------------------------------------
# Row 2 (helper text)
# ---------------------------------------------------------------------------
# B2:G2 are unchanged. Add the two new helper notes and an (empty, styled)
# placeholder cell under the new "ontological definition" column.
$ws.Range("J2").Value = "Parent term in ClinEpi (which genus in which sample type)"
$ws.Range("J2").WrapText = $true

$ws.Range("N2").WrapText = $true

# ---------------------------------------------------------------------------
# Row 3 (example data row)
# ---------------------------------------------------------------------------
# Sample type example changes from "stool" to "blood".
$ws.Range("B3").Value = "blood"

$ws.Range("K2").Value = "(which domain in which sample type)"
$ws.Range("K2").WrapText = $true

# New lookup-style formulas for website parent/grandparent.
$ws.Range("J3").Formula = '=F3&" in "&B3'
$ws.Range("K3").Formula = '=E3&" in "&B3'
$ws.Range("J3:K3").WrapText = $true

# The old "ontological label"/"ontological definition" formulas move from
# K3/L3 to M3/N3.
$ws.Range("M3").Formula = '=IF(D3="boolean","presence of",IF(D3="count","count of","data about"))&" "&G3&" by "&C3'
$ws.Range("N3").Formula = '=IF($D3="count","a count of the number of ",IF($D3="boolean","a categorical measurement datum","a data item")&" that is about ")&$G3&" and is the specified output of some "&C3&" assay, which achieves an organism identification objective and has as specified input a "&B3&" specimen from an organism"'
$ws.Range("M3:N3").WrapText = $true

# L3 no longer holds data (its formula moved to N3) - remove it entirely.
$ws.Range("L3").Clear()

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
# Column J widens from its old narrow width to match column K.
$ws.Range("J1").ColumnWidth = $ws.Range("K1").ColumnWidth
# New column N gets its own width.
$ws.Range("N1").ColumnWidth = 32

# ---------------------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------------------
$ws.Range("A1").RowHeight = 51

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("N3").Select()
